$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B23 should become a numeric value 3 (was an inline string "3")
$ws.Range("B23").Value = 3

# Add new row 24 with the "Cons" annotation entry
$ws.Range("A24").Value = "Ruilin"
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = "Cons: - it does not seem to involve any learning, it clearly does not fit at ICLR."
$ws.Range("D24").Value = "CRT"
$ws.Range("E24").Value = "OTH"
$ws.Range("F24").Value = "b49eb73e-9ff0-45de-a177-7d78dc315c92"
$ws.Range("G24").Value = "2rHk2kZ5knTJ6_annotated.xlsx"
$ws.Range("H24").Value = "Cons: - it does not seem to involve any learning, it clearly does not fit at ICLR."
